$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.260.39"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "2.639.94"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.17"
$ws.Range("E5").Value = "  -2.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.74"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.24%  "

$ws.Range("D9").Value = "2.638.48"
$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("E13").Value = "  -2.23%  "

$ws.Range("D14").Value = "3.132.02"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("E15").Value = "  -1.85%  "

$ws.Range("D16").Value = "72.159.28"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.88"
$ws.Range("E17").Value = "  -2.62%  "

$ws.Range("D18").Value = "2.658.03"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.21"
$ws.Range("E19").Value = "  +1.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.97"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.70"
$ws.Range("E21").Value = "  -1.49%  "

$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("E23").Value = "  -1.68%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.89"
$ws.Range("E25").Value = "  -2.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.27"
$ws.Range("E26").Value = "  -3.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  -3.86%  "

$ws.Range("D28").Value = "2.774.87"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "0.0₃0956"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("E31").Value = "  -2.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "497.02"
$ws.Range("E32").Value = "  -4.68%  "

$ws.Range("E33").Value = "  -3.17%  "

$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.02"
$ws.Range("E36").Value = "  -1.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.29"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.113"
$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.88"
$ws.Range("E39").Value = "  -1.23%  "

$ws.Range("E40").Value = "  -2.99%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  -6.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.56"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.90"
$ws.Range("E44").Value = "  -4.09%  "

$ws.Range("E45").Value = "  -2.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.09"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.12"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.66"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.548"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("E51").Value = "  -1.87%  "
